$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 values
$ws.Range("B2").Value = "E:\storage\114G.mp4"
$ws.Range("C2").Value = "đâssad"
$ws.Range("G2").Value = "E:/New folder\114G.mp4"

# Clear D2, E2, F2 (no longer present in the data)
$ws.Range("D2").ClearContents()
$ws.Range("E2").ClearContents()
$ws.Range("F2").ClearContents()

# Delete rows 3 through 5 entirely
$ws.Range("A3:G5").Delete()
